$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 15627266
$ws.Range("I132").Value = 18520242
$ws.Range("J132").Value = 5191.7
$ws.Range("K132").Value = 55560726
$ws.Range("L132").Value = 15575.1
$ws.Range("M132").Value = -55558196
$ws.Range("N132").Value = -20635.1
$ws.Range("H137").Value = 5017.8125
$ws.Range("I137").Value = 7072.1816
$ws.Range("J137").Value = 3941.7144
$ws.Range("K137").Value = 21216.5448
$ws.Range("L137").Value = 11825.1432
$ws.Range("M137").Value = -18666.5448
$ws.Range("N137").Value = -16925.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1949.909
$ws.Range("I32").Value = 1736.4445
$ws.Range("J32").Value = 2910.5
$ws.Range("K32").Value = 1736.4445
$ws.Range("L32").Value = 2910.5
$ws.Range("M32").Value = -1449.4445
$ws.Range("N32").Value = -3484.5
$ws.Range("H61").Value = 1535.6786
$ws.Range("I61").Value = 1140.2273
$ws.Range("J61").Value = 2985.6667
$ws.Range("K61").Value = 1140.2273
$ws.Range("L61").Value = 2985.6667
$ws.Range("M61").Value = -928.2273
$ws.Range("N61").Value = -3409.6667
$ws.Range("H74").Value = 2979.9487
$ws.Range("I74").Value = 3431.8518
$ws.Range("J74").Value = 1963.1666
$ws.Range("K74").Value = 3431.8518
$ws.Range("L74").Value = 1963.1666
$ws.Range("M74").Value = -2557.8518
$ws.Range("N74").Value = -3711.1666
$ws.Range("H77").Value = 2979.9487
$ws.Range("I77").Value = 3431.8518
$ws.Range("J77").Value = 1963.1666
$ws.Range("K77").Value = 17159.259
$ws.Range("L77").Value = 9815.833000000001
$ws.Range("M77").Value = -12791.259
$ws.Range("N77").Value = -18551.833
$ws.Range("H110").Value = 1698.1111
$ws.Range("I110").Value = 1213.8334
$ws.Range("J110").Value = 2666.6667
$ws.Range("K110").Value = 1213.8334
$ws.Range("L110").Value = 2666.6667
$ws.Range("M110").Value = 831.1666
$ws.Range("N110").Value = -6756.6667
$ws.Range("H132").Value = 1606.638
$ws.Range("I132").Value = 995.8431399999999
$ws.Range("J132").Value = 6056.7144
$ws.Range("K132").Value = 2987.52942
$ws.Range("L132").Value = 18170.1432
$ws.Range("M132").Value = -457.5294199999998
$ws.Range("N132").Value = -23230.1432
$ws.Range("H136").Value = 1535.6786
$ws.Range("I136").Value = 1140.2273
$ws.Range("J136").Value = 2985.6667
$ws.Range("K136").Value = 3420.6819
$ws.Range("L136").Value = 8957.000100000001
$ws.Range("M136").Value = -870.6819
$ws.Range("N136").Value = -14057.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 559.8
$ws.Range("I64").Value = 549.5
$ws.Range("K64").Value = 549.5
$ws.Range("M64").Value = -324.5
$ws.Range("H67").Value = 559.8
$ws.Range("I67").Value = 549.5
$ws.Range("K67").Value = 549.5
$ws.Range("M67").Value = 230.5
$ws.Range("H134").Value = 2138.0833
$ws.Range("I134").Value = 1258.4667
$ws.Range("J134").Value = 4776.933
$ws.Range("K134").Value = 3775.4001
$ws.Range("L134").Value = 14330.799
$ws.Range("M134").Value = -1240.4001
$ws.Range("N134").Value = -19400.799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12503373
$ws.Range("I31").Value = 1775.421
$ws.Range("J31").Value = 23814342
$ws.Range("K31").Value = 1775.421
$ws.Range("L31").Value = 23814342
$ws.Range("M31").Value = -1480.421
$ws.Range("N31").Value = -23814932
$ws.Range("H34").Value = 12503373
$ws.Range("I34").Value = 1775.421
$ws.Range("J34").Value = 23814342
$ws.Range("K34").Value = 1775.421
$ws.Range("L34").Value = 23814342
$ws.Range("M34").Value = -1573.421
$ws.Range("N34").Value = -23814746
$ws.Range("H99").Value = 33340996
$ws.Range("I99").Value = 100003000
$ws.Range("J99").Value = 9995
$ws.Range("K99").Value = 100003000
$ws.Range("L99").Value = 9995
$ws.Range("M99").Value = -100001502
$ws.Range("N99").Value = -12991
$ws.Range("H126").Value = 33340996
$ws.Range("I126").Value = 100003000
$ws.Range("J126").Value = 9995
$ws.Range("K126").Value = 300009000
$ws.Range("L126").Value = 29985
$ws.Range("M126").Value = -300006530
$ws.Range("N126").Value = -34925
$ws.Range("H132").Value = 1639.9436
$ws.Range("I132").Value = 1363.2097
$ws.Range("J132").Value = 3546.3333
$ws.Range("K132").Value = 4089.6291
$ws.Range("L132").Value = 10638.9999
$ws.Range("M132").Value = -1559.6291
$ws.Range("N132").Value = -15698.9999
$ws.Range("H134").Value = 1241.4324
$ws.Range("I134").Value = 746.8302
$ws.Range("J134").Value = 2489.7144
$ws.Range("K134").Value = 2240.4906
$ws.Range("L134").Value = 7469.1432
$ws.Range("M134").Value = 294.5093999999999
$ws.Range("N134").Value = -12539.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 624.1852
$ws.Range("I113").Value = 521.0909
$ws.Range("J113").Value = 1077.8
$ws.Range("K113").Value = 1563.2727
$ws.Range("L113").Value = 3233.4
$ws.Range("M113").Value = 606.7273
$ws.Range("N113").Value = -7573.4
$ws.Range("H131").Value = 758.4253
$ws.Range("I131").Value = 429.16666
$ws.Range("J131").Value = 811.1067
$ws.Range("K131").Value = 1287.49998
$ws.Range("L131").Value = 2433.3201
$ws.Range("M131").Value = 3752.50002
$ws.Range("N131").Value = -12513.3201

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1582.3243
$ws.Range("I102").Value = 1065.1034
$ws.Range("J102").Value = 3457.25
$ws.Range("K102").Value = 1065.1034
$ws.Range("L102").Value = 3457.25
$ws.Range("M102").Value = 556.8966
$ws.Range("N102").Value = -6701.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2048.3809
$ws.Range("J93").Value = 2379.4167
$ws.Range("L93").Value = 2379.4167
$ws.Range("N93").Value = -4875.4167
$ws.Range("H122").Value = 7004.909
$ws.Range("I122").Value = 2763.5
$ws.Range("J122").Value = 9428.571
$ws.Range("K122").Value = 8290.5
$ws.Range("L122").Value = 28285.713
$ws.Range("M122").Value = -5840.5
$ws.Range("N122").Value = -33185.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3159.0278
$ws.Range("I136").Value = 1190.7778
$ws.Range("J136").Value = 5127.278
$ws.Range("K136").Value = 3572.3334
$ws.Range("L136").Value = 15381.834
$ws.Range("M136").Value = -1022.3334
$ws.Range("N136").Value = -20481.834
